$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.218
$ws.Range("A7").Value = -20.064
$ws.Range("A16").Value = -22.119
$ws.Range("A28").Value = -22.005
$ws.Range("A29").Value = -21.344
$ws.Range("A32").Value = -21.782
$ws.Range("A40").Value = -19.965
$ws.Range("A52").Value = -21.957
$ws.Range("A57").Value = -22.253
$ws.Range("A66").Value = -21.53
$ws.Range("A100").Value = -22.352
